# Commit: "Allowing resource optimized OpenStack"
#
# The only substantive content change in the target revision is on slide 1
# (sldId 290), where the small vertical "OpenStack" label box loses its
# "LW " prefix, so the same box can stand for any (including resource
# optimized) OpenStack flavour instead of only the "LW" one.
#
# That shape is nested inside the big top-level group "组合 30" on the
# slide, so we recurse through GroupItems to find it by its stable shape
# Id (16 / creationId {F8ED39A8-1377-4641-955B-B363CC431355}) rather than
# relying on a fragile flat shape index.

function Find-ShapeById($shapes, $targetId) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Id -eq $targetId) {
            return $shp
        }
        if ($shp.Type -eq 6) {
            # msoGroup - recurse into the group's members
            $found = Find-ShapeById $shp.GroupItems $targetId
            if ($found -ne $null) {
                return $found
            }
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)

$openStackShape = Find-ShapeById $slide.Shapes 16

if ($openStackShape -ne $null -and $openStackShape.HasTextFrame) {
    $openStackShape.TextFrame.TextRange.Text = "OpenStack"
}
